$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Fitness), rows 2 through 252 all change from 7573 to 7293
$ws.Range("C2:C252").Value = 7293
